$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Capture reference fill colors before touching anything, so style changes
# --- reuse the workbook's existing style slots instead of minting new ones.
$greenColor = $ws1.Range("C3").Interior.Color
$redColor   = $ws1.Range("G3").Interior.Color

# --- Text fields (A3 "name", H3 "tempo de prova") are plain strings, safe to
# --- set directly via .Value on both sheets.
$ws1.Range("A3").Value = "OH KAMI OH MY GOD"
$ws2.Range("A3").Value = "OH KAMI OH MY GOD"

$ws1.Range("H3").Value = "00:43:20"
$ws2.Range("H3").Value = "00:43:20"

# --- B3 "date" text looks like a real date (12/12/2012), so a plain .Value
# --- assignment gets auto-converted to a date serial + new number format.
# --- Write it as text in a scratch cell, then paste-values into B3 so the
# --- destination keeps its original (unchanged) cell style.
$ws1.Range("ZZ1").NumberFormat = "@"
$ws1.Range("ZZ1").Value = "12/12/2012"
$ws1.Range("ZZ1").Copy() | Out-Null
$ws1.Range("B3").PasteSpecial(-4163) | Out-Null
$ws2.Range("B3").PasteSpecial(-4163) | Out-Null
$ws1.Range("ZZ1").Clear() | Out-Null
$excel.CutCopyMode = 0

# --- Prova1 (sheet1) row 3 numeric cells + fill-color swaps
$ws1.Range("D3").Value = 1000
$ws1.Range("D3").Interior.Color = $greenColor

$ws1.Range("E3").Value = 750
$ws1.Range("E3").Interior.Color = $greenColor

$ws1.Range("F3").Value = 300
$ws1.Range("F3").Interior.Color = $redColor

$ws1.Range("G3").Value = 200

$ws1.Range("I3").Value = 0

# --- Prova2 (sheet2) row 3 numeric cells (fill colors unchanged)
$ws2.Range("C3").Value = 0
$ws2.Range("D3").Value = 676
$ws2.Range("E3").Value = 600
$ws2.Range("F3").Value = 80
$ws2.Range("G3").Value = 40
$ws2.Range("I3").Value = 0
